# Fix mo hinh As, update banner
# Appends new "ngayxem" (viewing-date) log rows to the dataset sheet,
# reusing the existing shared-string date values where applicable and
# introducing five new date strings (2024-06-23, 2024-06-26, 2024-06-27,
# 2024-06-28, 2024-06-29) for the newly logged days.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of the new "ngayxem" shared-string indices (as found in the
# post-edit sharedStrings.xml) to their date text.
$dateByIndex = @{
    20 = "2024-06-23"
    21 = "2024-06-26"
    22 = "2024-06-27"
    23 = "2024-06-28"
    24 = "2024-06-29"
}

# Each entry: row number, date-index (see $dateByIndex), masp (col B), makh (col C)
$newRows = @(
    @(794, 20, 49, 62),
    @(795, 21, 49, 65),
    @(796, 21, 49, 65),
    @(797, 21, 49, 65),
    @(798, 21, 49, 65),
    @(799, 21, 49, 65),
    @(800, 21, 49, 65),
    @(801, 21, 49, 65),
    @(802, 21, 49, 65),
    @(803, 21, 49, 65),
    @(804, 21, 49, 65),
    @(805, 21, 49, 65),
    @(806, 21, 49, 65),
    @(807, 21, 49, 65),
    @(808, 21, 49, 65),
    @(809, 21, 49, 65),
    @(810, 21, 49, 65),
    @(811, 21, 49, 65),
    @(812, 21, 49, 65),
    @(813, 21, 49, 65),
    @(814, 21, 49, 65),
    @(815, 21, 49, 65),
    @(816, 21, 49, 65),
    @(817, 21, 49, 65),
    @(818, 21, 49, 65),
    @(819, 21, 49, 65),
    @(820, 21, 49, 65),
    @(821, 21, 49, 65),
    @(822, 21, 49, 65),
    @(823, 21, 49, 65),
    @(824, 21, 49, 65),
    @(825, 21, 49, 65),
    @(826, 21, 49, 65),
    @(827, 21, 49, 65),
    @(828, 21, 49, 65),
    @(829, 21, 49, 65),
    @(830, 21, 49, 65),
    @(831, 21, 49, 65),
    @(832, 21, 49, 65),
    @(833, 21, 77, 65),
    @(834, 21, 77, 65),
    @(835, 21, 76, 61),
    @(836, 22, 49, 66),
    @(837, 22, 103, 52),
    @(838, 22, 103, 52),
    @(839, 22, 103, 52),
    @(840, 22, 86, 52),
    @(841, 22, 58, 52),
    @(842, 22, 64, 52),
    @(843, 22, 49, 52),
    @(844, 22, 144, 52),
    @(845, 23, 50, 52),
    @(846, 23, 103, 53),
    @(847, 23, 72, 53),
    @(848, 23, 56, 53),
    @(849, 23, 145, 53),
    @(850, 23, 139, 53),
    @(851, 23, 149, 53),
    @(852, 23, 153, 53),
    @(853, 23, 149, 53),
    @(854, 23, 144, 33),
    @(855, 23, 145, 33),
    @(856, 23, 73, 33),
    @(857, 23, 136, 33),
    @(858, 23, 145, 33),
    @(859, 23, 49, 67),
    @(860, 23, 49, 67),
    @(861, 23, 49, 67),
    @(862, 23, 49, 67),
    @(863, 23, 136, 67),
    @(864, 23, 71, 67),
    @(865, 23, 182, 67),
    @(866, 23, 50, 67),
    @(867, 23, 215, 67),
    @(868, 23, 144, 28),
    @(869, 23, 64, 27),
    @(870, 23, 145, 61),
    @(871, 24, 172, 28),
    @(872, 24, 169, 28)
)

foreach ($entry in $newRows) {
    $rowNum = $entry[0]
    $dateIdx = $entry[1]
    $masp = $entry[2]
    $makh = $entry[3]

    # Write column A as literal text (not an auto-converted date serial) by
    # writing it as a formula that evaluates to the literal string, then
    # collapsing the formula to its static value via copy/paste-values.
    # This keeps the cell a plain shared-string text cell with the default
    # style, matching the rest of the "ngayxem" column.
    $cellA = $ws.Range("A$rowNum")
    $cellA.Formula = '="' + $dateByIndex[$dateIdx] + '"'
    $cellA.Copy()
    $cellA.PasteSpecial(-4163)

    $ws.Range("B$rowNum").Value = $masp
    $ws.Range("C$rowNum").Value = $makh
}
